# Added 4wk low sales check
# Updates forecast numbers and derived risk/urgency metrics on the
# "Forecast Comparison" sheet, plus the roll-up totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("D2").Value = 2
$ws1.Range("L2").Value = 0.98

$ws1.Range("D3").Value = 2
$ws1.Range("L3").Value = 0.98

$ws1.Range("L4").Value = 0.88

$ws1.Range("L5").Value = 0.85

$ws1.Range("L6").Value = 1.11

$ws1.Range("L7").Value = 1

$ws1.Range("L8").Value = 0.93

$ws1.Range("L9").Value = 0.91

$ws1.Range("L10").Value = 0.93

$ws1.Range("L11").Value = 1.12

$ws1.Range("L12").Value = 1

$ws1.Range("L13").Value = 1.17

$ws1.Range("L14").Value = 0.96

$ws1.Range("L15").Value = 1.03

$ws1.Range("H16").Value = 0
$ws1.Range("I16").Value = "High"
$ws1.Range("J16").Value = "Urgent"
$ws1.Range("L16").Value = 0.82

$ws1.Range("D17").Value = 1
$ws1.Range("H17").Value = 0
$ws1.Range("I17").Value = "High"
$ws1.Range("J17").Value = "Urgent"
$ws1.Range("L17").Value = 1.18

# --- Sheet 2: Summary ---
# Column B cells hold numbers stored as text; use a leading apostrophe so
# Excel keeps them as text (matching the original "numeric-looking string"
# cell type) instead of converting them to real numbers.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "'16"
$ws2.Range("B10").Value = "'9"
$ws2.Range("B11").Value = "'6"
$ws2.Range("B12").Value = "'2"
$ws2.Range("B14").Value = "'1"
